# Updates cryptos list values/percentages (and one coin swap at row 51)
# Forces text-typed cell values (matching the source t="inlineStr" cells)
# by toggling NumberFormat to Text before the write, then resetting the
# cell style back to "Normal" so no stray style index is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "70.868.78"
Set-TextValue "E2" "  +6.96%  "
Set-TextValue "D3" "3.628.65"
Set-TextValue "E3" "  +6.12%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "592.33"
Set-TextValue "E5" "  +4.37%  "
Set-TextValue "D6" "192.04"
Set-TextValue "E6" "  +7.46%  "
Set-TextValue "E7" "  +3.10%  "
Set-TextValue "D8" "3.623.03"
Set-TextValue "E8" "  +6.19%  "
Set-TextValue "D9" "1.00"
Set-TextValue "E9" "  -0.02%  "
Set-TextValue "E11" "  +4.09%  "
Set-TextValue "D12" "58.35"
Set-TextValue "E12" "  +6.82%  "
Set-TextValue "E13" "  +6.24%  "
Set-TextValue "D14" "9.81"
Set-TextValue "E14" "  +5.24%  "
Set-TextValue "D15" "4.206.34"
Set-TextValue "E15" "  +6.25%  "
Set-TextValue "D16" "3.629.63"
Set-TextValue "E16" "  +6.09%  "
Set-TextValue "D17" "19.41"
Set-TextValue "D18" "70.779.51"
Set-TextValue "E18" "  +7.03%  "
Set-TextValue "D19" "12.58"
Set-TextValue "E19" "  +4.83%  "
Set-TextValue "E20" "  +0.47%  "
Set-TextValue "E21" "  +4.45%  "
Set-TextValue "D22" "496.69"
Set-TextValue "E22" "  +6.77%  "
Set-TextValue "D23" "17.32"
Set-TextValue "E23" "  +18.36%  "
Set-TextValue "D24" "5.40"
Set-TextValue "E24" "  +8.66%  "
Set-TextValue "D25" "4.51"
Set-TextValue "E25" "  +8.88%  "
Set-TextValue "D26" "91.11"
Set-TextValue "E26" "  +1.37%  "
Set-TextValue "E27" "  +6.53%  "
Set-TextValue "D28" "11.27"
Set-TextValue "E28" "  +4.67%  "
Set-TextValue "E29" "  +6.86%  "
Set-TextValue "D30" "32.45"
Set-TextValue "E30" "  +3.42%  "
Set-TextValue "D31" "7.57"
Set-TextValue "E31" "  +11.57%  "
Set-TextValue "D32" "12.26"
Set-TextValue "E32" "  +5.84%  "
Set-TextValue "D33" "622.88"
Set-TextValue "E33" "  +7.28%  "
Set-TextValue "D35" "65.27"
Set-TextValue "E35" "  +4.25%  "
Set-TextValue "D36" "0.0₃0833"
Set-TextValue "E36" "  +8.18%  "
Set-TextValue "D37" "0.412"
Set-TextValue "E37" "  +7.57%  "
Set-TextValue "D38" "38.26"
Set-TextValue "E38" "  +5.01%  "
Set-TextValue "E39" "  +1.73%  "
Set-TextValue "E40" "  +0.09%  "
Set-TextValue "E41" "  +2.48%  "
Set-TextValue "D42" "3.327.75"
Set-TextValue "E42" "  +6.31%  "
Set-TextValue "E43" "  +4.55%  "
Set-TextValue "D44" "0.0449"
Set-TextValue "E44" "  +6.03%  "
Set-TextValue "E45" "  +7.76%  "
Set-TextValue "D46" "3.34"
Set-TextValue "E46" "  +4.69%  "
Set-TextValue "E47" "  +2.77%  "
Set-TextValue "E48" "  +7.26%  "
Set-TextValue "E49" "  +2.90%  "
Set-TextValue "E50" "  +4.53%  "
Set-TextValue "B51" "Monero"
Set-TextValue "C51" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D51" "142.67"
Set-TextValue "E51" "  +0.72%  "
